$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.416200000000004
$ws.Range("C3").Value = -11.95970000000001
$ws.Range("D8").Value = -9.0725
$ws.Range("D11").Value = -7.309100000000001
$ws.Range("A12").Value = -21.61810000000001
$ws.Range("B14").Value = 5.9338
$ws.Range("D14").Value = -7.691999999999999
$ws.Range("D15").Value = -8.373599999999994
$ws.Range("D17").Value = -8.545699999999998
$ws.Range("C20").Value = -11.8787
$ws.Range("C25").Value = -12.4477
$ws.Range("B26").Value = 4.207700000000006
$ws.Range("D26").Value = -8.612700000000007
$ws.Range("A27").Value = -21.66159999999999
$ws.Range("C30").Value = -13.54659999999999
$ws.Range("B31").Value = 4.792600000000002
$ws.Range("A32").Value = -21.48730000000002
$ws.Range("B35").Value = 9.067400000000006
$ws.Range("A36").Value = -19.6855
$ws.Range("D36").Value = -7.586999999999998
$ws.Range("B37").Value = 9.144199999999998
$ws.Range("A38").Value = -19.45420000000001
$ws.Range("C44").Value = -13.54889999999999
$ws.Range("B45").Value = 5.7323
$ws.Range("A46").Value = -21.71240000000001
$ws.Range("C47").Value = -12.23039999999999
$ws.Range("B52").Value = 5.436700000000001
$ws.Range("A54").Value = -21.71849999999998
$ws.Range("A55").Value = -22.1969
$ws.Range("A56").Value = -22.10510000000001
$ws.Range("B57").Value = 4.825199999999996
$ws.Range("C58").Value = -13.3273
$ws.Range("D64").Value = -7.437699999999998
$ws.Range("A67").Value = -21.50649999999998
$ws.Range("A69").Value = -21.62059999999998
$ws.Range("A72").Value = -21.85709999999999
$ws.Range("C78").Value = -10.91100000000001
$ws.Range("D79").Value = -6.231700000000001
$ws.Range("B81").Value = 6.2043
$ws.Range("A83").Value = -21.6796
$ws.Range("B83").Value = 5.727700000000002
$ws.Range("C84").Value = -13.91329999999999
$ws.Range("A86").Value = -22.087
$ws.Range("C89").Value = -11.4061
$ws.Range("D89").Value = -6.109699999999999
$ws.Range("A91").Value = -21.5823
$ws.Range("C91").Value = -11.4796
$ws.Range("C92").Value = -11.7221
$ws.Range("A93").Value = -21.2944
$ws.Range("C96").Value = -13.8018
$ws.Range("A99").Value = -20.36419999999998
$ws.Range("B100").Value = 5.625099999999998
$ws.Range("B102").Value = 8.170100000000005
$ws.Range("C102").Value = -13.7496
